# Updates the "Price" (column D) and "Volume(1h)" (column E) figures
# in the cryptocurrency price table on the active sheet, matching the
# refreshed data pulled by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (column D, omitted if unchanged),
# whether the Price string parses as a plain number (needs to be forced
# to Text so Excel keeps it verbatim instead of reformatting it), and the
# new Volume(1h) percentage (column E, always updated).
$updates = @(
    @{ Row = 2; Coin = "Bitcoin"; Price = "64.369.15"; PriceIsNumber = $false; Volume = "+0.11%" }
    @{ Row = 3; Coin = "Ethereum"; Price = "3.513.92"; PriceIsNumber = $false; Volume = "+0.28%" }
    @{ Row = 4; Coin = "TetherUSD"; Price = $null; PriceIsNumber = $false; Volume = "+0.04%" }
    @{ Row = 5; Coin = "BNB"; Price = "591.60"; PriceIsNumber = $true; Volume = "+1.17%" }
    @{ Row = 6; Coin = "Solana"; Price = "134.65"; PriceIsNumber = $true; Volume = "-0.17%" }
    @{ Row = 8; Coin = "XRP"; Price = $null; PriceIsNumber = $false; Volume = "+0.06%" }
    @{ Row = 9; Coin = "Toncoin"; Price = "7.63"; PriceIsNumber = $true; Volume = "+7.15%" }
    @{ Row = 10; Coin = "Dogecoin"; Price = $null; PriceIsNumber = $false; Volume = "+0.25%" }
    @{ Row = 11; Coin = "Cardano"; Price = $null; PriceIsNumber = $false; Volume = "+4.03%" }
    @{ Row = 12; Coin = "WrappedliquidstakedEther2.0"; Price = "4.114.47"; PriceIsNumber = $false; Volume = "+0.33%" }
    @{ Row = 14; Coin = "ShibaInu"; Price = "0.0000181"; PriceIsNumber = $true; Volume = "+0.98%" }
    @{ Row = 15; Coin = "WrappedEther"; Price = "3.513.78"; PriceIsNumber = $false; Volume = "+0.22%" }
    @{ Row = 16; Coin = "Avalanche"; Price = "25.86"; PriceIsNumber = $true; Volume = "-1.73%" }
    @{ Row = 17; Coin = "WrappedBTC"; Price = "64.355.17"; PriceIsNumber = $false; Volume = "+0.11%" }
    @{ Row = 18; Coin = "Uniswap"; Price = "10.00"; PriceIsNumber = $true; Volume = "+2.34%" }
    @{ Row = 19; Coin = "Polkadot"; Price = $null; PriceIsNumber = $false; Volume = "+3.24%" }
    @{ Row = 20; Coin = "Chainlink"; Price = "13.56"; PriceIsNumber = $true; Volume = "-2.13%" }
    @{ Row = 21; Coin = "BitcoinCash"; Price = "394.11"; PriceIsNumber = $true; Volume = "+2.73%" }
    @{ Row = 22; Coin = "Polygon"; Price = $null; PriceIsNumber = $false; Volume = "+1.10%" }
    @{ Row = 23; Coin = "WrappedeETH"; Price = "3.655.07"; PriceIsNumber = $false; Volume = "+0.38%" }
    @{ Row = 24; Coin = "Litecoin"; Price = "74.67"; PriceIsNumber = $true; Volume = "+0.89%" }
    @{ Row = 25; Coin = "Dai"; Price = $null; PriceIsNumber = $false; Volume = "+0.06%" }
    @{ Row = 27; Coin = "PEPE"; Price = $null; PriceIsNumber = $false; Volume = "+3.13%" }
    @{ Row = 28; Coin = "Binance-PegBSC-USD"; Price = $null; PriceIsNumber = $false; Volume = "+0.01%" }
    @{ Row = 31; Coin = "InternetComputer(DFINITY)"; Price = "8.30"; PriceIsNumber = $true; Volume = "+0.11%" }
    @{ Row = 32; Coin = "Fetch.AI"; Price = $null; PriceIsNumber = $false; Volume = "-6.52%" }
    @{ Row = 33; Coin = "Kaspa"; Price = $null; PriceIsNumber = $false; Volume = "+7.36%" }
    @{ Row = 34; Coin = "RenzoRestakedETH"; Price = "3.545.64"; PriceIsNumber = $false; Volume = "+0.63%" }
    @{ Row = 35; Coin = "USDe"; Price = $null; PriceIsNumber = $false; Volume = "+0.04%" }
    @{ Row = 36; Coin = "EthereumClassic"; Price = "23.42"; PriceIsNumber = $true; Volume = "-0.60%" }
    @{ Row = 37; Coin = "NEARProtocol"; Price = "5.35"; PriceIsNumber = $true; Volume = "+0.85%" }
    @{ Row = 38; Coin = "Aptos"; Price = $null; PriceIsNumber = $false; Volume = "+1.57%" }
    @{ Row = 39; Coin = "ImmutableX"; Price = $null; PriceIsNumber = $false; Volume = "+0.90%" }
    @{ Row = 40; Coin = "Monero"; Price = "167.01"; PriceIsNumber = $true; Volume = "+1.62%" }
    @{ Row = 41; Coin = "Hedera"; Price = $null; PriceIsNumber = $false; Volume = "+0.83%" }
    @{ Row = 42; Coin = "Mantle"; Price = $null; PriceIsNumber = $false; Volume = "+0.40%" }
    @{ Row = 43; Coin = "EnergySwap"; Price = "25.53"; PriceIsNumber = $true; Volume = "-3.06%" }
    @{ Row = 44; Coin = "FirstDigitalUSD"; Price = $null; PriceIsNumber = $false; Volume = "+0.02%" }
    @{ Row = 45; Coin = "Filecoin"; Price = "4.45"; PriceIsNumber = $true; Volume = "+0.84%" }
    @{ Row = 46; Coin = "Stacks"; Price = $null; PriceIsNumber = $false; Volume = "+2.93%" }
    @{ Row = 47; Coin = "ONDO"; Price = $null; PriceIsNumber = $false; Volume = "-2.10%" }
    @{ Row = 48; Coin = "Cosmos"; Price = $null; PriceIsNumber = $false; Volume = "+0.71%" }
    @{ Row = 49; Coin = "Maker"; Price = "2.393.91"; PriceIsNumber = $false; Volume = "-3.53%" }
    @{ Row = 50; Coin = "SuiNetwork"; Price = "0.899"; PriceIsNumber = $true; Volume = "-2.26%" }
    @{ Row = 51; Coin = "VeChain"; Price = $null; PriceIsNumber = $false; Volume = "+0.41%" }
)

foreach ($u in $updates) {
    if ($null -ne $u.Price) {
        $priceCell = $ws.Range("D" + $u.Row)
        if ($u.PriceIsNumber) {
            # Looks like a plain number (e.g. "591.60") - force Text format first
            # so Excel keeps the exact original string instead of converting it
            # to a numeric value (which would drop trailing zeros, etc.), then
            # restore the default "Normal" style so no stray formatting is left
            # behind on the cell.
            $priceCell.NumberFormat = "@"
            $priceCell.Value = $u.Price
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $u.Price
        }
    }
    $ws.Range("E" + $u.Row).Value = "  " + $u.Volume + "  "
}
